# Hortaliza, Terminal La Palmera de La Serena - Cebollín
# Insert a new weekly observation row at the top of the data block (row 279),
# pushing the existing rows 279-373 down to 280-374.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before row 279 (shifts 279..373 -> 280..374)
$ws.Rows.Item(279).EntireRow.Insert()

# Populate the new row 279 with the latest weekly price record
$ws.Cells.Item(279, 1).Value = 8
$ws.Cells.Item(279, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(279, 3).Value = "Coquimbo"
$ws.Cells.Item(279, 4).Value = 45215
$ws.Cells.Item(279, 5).Value = 4
$ws.Cells.Item(279, 6).Value = 100112037
$ws.Cells.Item(279, 7).Value = "Cebollín"
$ws.Cells.Item(279, 8).Value = "Sin especificar"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 1400
$ws.Cells.Item(279, 11).Value = 1000
$ws.Cells.Item(279, 12).Value = 1200
$ws.Cells.Item(279, 13).Value = 1100
$ws.Cells.Item(279, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(279, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(279, 16).Value = 183
$ws.Cells.Item(279, 17).Value = 6
$ws.Cells.Item(279, 18).Value = "Hortaliza"
